$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.984.14"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "1.741.48"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'246.99"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.90%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "'0.5046"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -4.64%  "
$ws.Range("D8").Value = "'0.2740"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.85%  "
$ws.Range("D9").Value = "'0.06176"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").Value = "1.754.35"
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("D11").Value = "'0.07247"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("D12").Value = "'0.6532"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("D13").Value = "'15.10"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("D14").Value = "'4.678"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("D15").Value = "'77.58"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").Value = "'0.9999"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "'1.000"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "26.004.46"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("D20").Value = "'0.000006855"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.36%  "
$ws.Range("D21").Value = "1.976.23"
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("D22").Value = "'4.474"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.24%  "
$ws.Range("D23").Value = "'8.718"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("D24").Value = "'5.369"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.88%  "
$ws.Range("D25").Value = "'135.80"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.25%  "
$ws.Range("D26").Value = "'1.516"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'15.25"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").Value = "'1.781"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("D29").Value = "'105.45"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("D30").Value = "'3.911"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.08%  "
$ws.Range("D31").Value = "'0.08165"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.98%  "
$ws.Range("E32").Value = "  +0.70%  "
$ws.Range("D33").Value = "'0.04677"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.76%  "
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").Value = "'0.9962"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("D36").Value = "'0.6123"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.07%  "
$ws.Range("D37").Value = "'2.763"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.43%  "
$ws.Range("D38").Value = "'0.01621"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("D39").Value = "'1.927"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.32%  "
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("D41").Value = "'100.79"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("D42").Value = "'0.3921"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.74%  "
$ws.Range("D43").Value = "'0.7639"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.50%  "
$ws.Range("D44").Value = "'5.006"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.05%  "
$ws.Range("D45").Value = "'0.1157"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.13%  "
$ws.Range("D46").Value = "'6.316"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.33%  "
$ws.Range("D47").Value = "'55.47"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("D48").Value = "'0.05301"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.19%  "
$ws.Range("D49").Value = "'30.67"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("D50").Value = "'0.3466"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("D51").Value = "'7.602"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.17%  "
